$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D40").Value = "[name=""마리아""]  아, 아니야……`n"
$ws.Range("D74").Value = "[name=""플래티넘""]  네가 걱정할 필욘 없어. 확실한 건, 이 모든 게 상업연합회의 명령이라는 거야.`n"
$ws.Range("D90").Value = "[name=""조피아""]  잉그라에 토플라까지, 사지 멀쩡하고 중상 안 입은 것만으로도 행운인 줄 알아 이것아……`n"
$ws.Range("D114").Value = "[name=""마리아""]  언니가 가장 싫어했던 게 기사 스포츠였을 거라고…… 사실 나도 한동안은 언니나 무에나 삼촌이나 생각이 아주 비슷하다고 생각했었거든……`n"
$ws.Range("D119").Value = "[name=""마리아""]  에헤헷…… 사실은 나도……`n"
$ws.Range("D125").Value = "[name=""조피아""]  왜, 왜 그래? 얘는! 왜 사람을 놀라게 하고 그래……`n"
